# Applies the "Added items into mapping, updated arduino to be compatible
# with NTags" commit:
#   - replace the 5 placeholder RFIDs (AAAAAAAA..EEEEEEEE) on the
#     rfid_item sheet/table with real NTag RFID codes
#   - give several items multiple RFID tags (new rows appended to Table2)
#   - tidy up the view state (active sheet/selection/zoom) to match
#     the saved session

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # device_cart
$ws2 = $wb.Worksheets.Item(2)   # rfid_item

# --- rfid_item: update existing rows with the new NTag RFID codes ----
$ws2.Range("A6").Value = "049B4B22"
$ws2.Range("A7").Value = "049F4B22"

$ws2.Range("A8").Value = "04974B22"
$ws2.Range("B8").Value = "VfgrHcX6LvHuAvkJtdgU"

$ws2.Range("A9").Value = "04934B22"
$ws2.Range("B9").Value = "YvxptylcQC7o6s7fK7H9"

$ws2.Range("A10").Value = "048F4B22"
$ws2.Range("B10").Value = "YvxptylcQC7o6s7fK7H9"

# --- rfid_item: append 4 new rows to Table2 for additional NTags -----
$lo = $ws2.ListObjects.Item("Table2")
$lo.ListRows.Add() | Out-Null
$lo.ListRows.Add() | Out-Null
$lo.ListRows.Add() | Out-Null
$lo.ListRows.Add() | Out-Null

$ws2.Range("A11").Value = "048B4B22"
$ws2.Range("B11").Value = "oZGiQLJMymfo2Mc4KJYm"

$ws2.Range("A12").Value = "04874B22"
$ws2.Range("B12").Value = "rxRod7cigQjBK9dDmlHv"

$ws2.Range("A13").Value = "04834B22"
$ws2.Range("B13").Value = "rxRod7cigQjBK9dDmlHv"

$ws2.Range("A14").Value = "047F4B22"
$ws2.Range("B14").Value = "rxRod7cigQjBK9dDmlHv"

# --- restore view state: rfid_item gets a zoom + new selection, -------
# --- but device_cart ends up the active/selected tab ------------------
$ws2.Activate()
$ws2.Range("E12").Select() | Out-Null
$excel.ActiveWindow.Zoom = 130

$ws1.Activate()
$ws1.Range("D15").Select() | Out-Null
